$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A8").Value = "Fixed lightmap in the nether(1.21.6)"
$ws.Range("D5").Select()
